# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Updates the "Metadata" sheet (canonical URL, version, status, date,
# description, context) and the "Elements" sheet (Extension definition text,
# Extension.value[x] cardinality Min 0 -> 1) of the UK Core Ethnic Category
# StructureDefinition export.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/UKCore-Extension-EthnicCategory"

# Version
$meta.Range("B3").Value = "0.1.0"

# Status
$meta.Range("B6").Value = "draft"

# Date
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description (was blank)
$meta.Range("B11").Value = "An extension to record the ethnic category of a patient, as per UK Core standards."

# Context
$meta.Range("B20").Value = "element:Patient"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" element: Definition column (M) text update
$elements.Range("M2").Value = "An extension to record the ethnic category of a patient, as per UK Core standards."

# Row 5 = "Extension.url" element: Fixed Value column (R) mirrors the
# canonical URL shown on the Metadata sheet, so it must be updated too.
$elements.Range("R5").Value = "https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/UKCore-Extension-EthnicCategory"

# Row 6 = "Extension.value[x]" element: Min column (F) 0 -> 1
# (kept as text, like the rest of the Min/Max columns in this table)
$elements.Range("F6").NumberFormat = "@"
$elements.Range("F6").Value = "1"
